$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("25-10-2021", $null, 2.74, 3.41),
    @("26-10-2021", $null, $null, 3.6),
    @("27-10-2021", $null, $null, 3.74),
    @("28-10-2021", $null, 2.85, 3.4),
    @("29-10-2021", $null, $null, 3.41)
)

$startRow = 189
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    if ($null -ne $rowData[2]) {
        $ws.Cells.Item($row, 3).Value = $rowData[2]
    }
    if ($null -ne $rowData[3]) {
        $ws.Cells.Item($row, 4).Value = $rowData[3]
    }
}
